# Weekly update: a new price record for Mandarina (Agrícola del Norte S.A. de Arica)
# is inserted as the new row 42, pushing all existing rows (old 42..125) down by one
# (new 43..126). This matches how the source sheet is maintained: newest record goes
# on top of this block, everything below shifts down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 42 - shifts rows 42:125 down to 43:126
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C42").Value = "Arica y Parinacota"
$ws.Range("D42").Value = 44868
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100102
$ws.Range("H42").Value = "Cítricos"
$ws.Range("I42").Value = 100102004
$ws.Range("J42").Value = "Mandarina"
$ws.Range("K42").Value = "Murcott"
$ws.Range("L42").Value = "Tercera"
$ws.Range("M42").Value = 350
$ws.Range("N42").Value = 13000
$ws.Range("O42").Value = 14000
$ws.Range("P42").Value = 13429
$ws.Range("Q42").Value = "$/caja 20 kilos"
$ws.Range("R42").Value = "Región de Coquimbo"
$ws.Range("S42").Value = 671
$ws.Range("T42").Value = 20
